$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (header row 3 / data row 4) gains a new trailing "2020" column
# (L) that duplicates the existing "2020" column (K) - both the year header
# and its data value. Copying K -> L brings along the same cell styling
# (borders/alignment/number format) that the rest of the yearly columns use.
$ws.Range("K3").Copy($ws.Range("L3")) | Out-Null
$ws.Range("K4").Copy($ws.Range("L4")) | Out-Null

# Reflect the new selection left behind after the edit.
$ws.Range("L10").Select() | Out-Null
